# Apply cryptos list price/volume refresh (GitHub Actions commit Fri Sep  1 19:29:39 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.026.90"
$ws.Range("E2").Value = "  -1.31%  "

$ws.Range("D3").Value = "1.639.75"
$ws.Range("E3").Value = "  -1.53%  "

$ws.Range("D4").Value = "'1.019"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.06%  "

$ws.Range("D5").Value = "'215.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "

$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("D7").Value = "'0.5004"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.87%  "

$ws.Range("D8").Value = "'0.2575"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("D9").Value = "'0.06421"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.29%  "

$ws.Range("D10").Value = "'19.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.33%  "

$ws.Range("D11").Value = "'0.07773"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.38%  "

$ws.Range("D12").Value = "1.647.69"
$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("D13").Value = "'4.257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.81%  "

$ws.Range("D14").Value = "1.864.42"
$ws.Range("E14").Value = "  -1.58%  "

$ws.Range("D15").Value = "'0.5450"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.45%  "

$ws.Range("D16").Value = "0.0₅7920"
$ws.Range("E16").Value = "  -1.30%  "

$ws.Range("D17").Value = "'63.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.45%  "

$ws.Range("D18").Value = "26.005.58"
$ws.Range("E18").Value = "  -1.52%  "

$ws.Range("D19").Value = "'1.018"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "

$ws.Range("D20").Value = "'203.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.65%  "

$ws.Range("D21").Value = "'4.310"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.31%  "

$ws.Range("D22").Value = "'9.995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.91%  "

$ws.Range("E23").Value = "  +1.71%  "

$ws.Range("E24").Value = "  +0.73%  "

$ws.Range("D25").Value = "'1.973"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +13.88%  "

$ws.Range("D26").Value = "'141.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.61%  "

$ws.Range("D27").Value = "'0.1152"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.94%  "

$ws.Range("D28").Value = "'15.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "'6.795"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.77%  "

$ws.Range("D30").Value = "'0.05040"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.53%  "

$ws.Range("D31").Value = "'1.242"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.53%  "

$ws.Range("D32").Value = "'3.265"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.01%  "

$ws.Range("D33").Value = "'3.201"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.37%  "

$ws.Range("D34").Value = "'1.542"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.81%  "

$ws.Range("D35").Value = "'2.355"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.84%  "

$ws.Range("D36").Value = "'0.8915"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.63%  "

$ws.Range("D37").Value = "'2.618"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.94%  "

$ws.Range("D38").Value = "'0.5647"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.07%  "

$ws.Range("D39").Value = "1.125.76"
$ws.Range("E39").Value = "  -2.23%  "

$ws.Range("D40").Value = "'0.01562"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.06%  "

$ws.Range("D41").Value = "'2.584"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.53%  "

$ws.Range("D42").Value = "'1.015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("D43").Value = "'5.636"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "'0.8165"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.87%  "

$ws.Range("D45").Value = "'99.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.38%  "

$ws.Range("D46").Value = "1.775.38"
$ws.Range("E46").Value = "  -1.59%  "

$ws.Range("E47").Value = "  +1.90%  "

$ws.Range("D48").Value = "'0.4565"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.60%  "

$ws.Range("D49").Value = "'1.016"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.20%  "

$ws.Range("D50").Value = "'54.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.87%  "

$ws.Range("D51").Value = "'0.05042"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.36%  "
